$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "(308051846, Eyal  Sofer: -3,3)"
$ws.Range("B1").Value = "(312049950, Molham  Peretz: 3,-6)"
$ws.Range("C1").Value = "(308073899, Anan  Kirshenbaum: 4,1)"
$ws.Range("D1").Value = "(318869187, Soaad  Leibovich: 1,-1)"
$ws.Range("E1").Value = "(205898513, Asaf  Braymok: 5,-7)"
$ws.Range("F1").Value = "(318428158, Tal  Asulin: 6,5)"
$ws.Range("G1").Value = "(316028364, Sami  Castro: -3,2)"

$ws.Range("A3").Value = "cost: 343.78382769538166"
$ws.Range("A4").Value = "time: 63.75676553907633"
